$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.584.35"
$ws.Range("E2").Value = "  +2.18%  "
$ws.Range("D3").Value = "2.288.67"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "157.80"
$ws.Range("E5").Value = "  +15,655.70%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "306.90"
$ws.Range("E6").Value = "  +1.29%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "96.81"
$ws.Range("E7").Value = "  +6.09%  "
$ws.Range("E8").Value = "  +0.63%  "
$ws.Range("E9").Value = "  +0.23%  "
$ws.Range("E10").Value = "  +3.92%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "36.51"
$ws.Range("E11").Value = "  +13.86%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0805"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("E14").Value = "  +2.76%  "
$ws.Range("D15").Value = "2.644.40"
$ws.Range("E15").Value = "  +1.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.59"
$ws.Range("E16").Value = "  +2.98%  "
$ws.Range("D17").Value = "2.311.19"
$ws.Range("E17").Value = "  +2.44%  "
$ws.Range("E18").Value = "  +6.22%  "
$ws.Range("D19").Value = "42.464.90"
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  +4.59%  "
$ws.Range("D21").Value = "0.0₃0920"
$ws.Range("E21").Value = "  +2.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.01"
$ws.Range("E22").Value = "  +2.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.94"
$ws.Range("E23").Value = "  +2.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "243.60"
$ws.Range("E24").Value = "  +1.39%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.61"
$ws.Range("E25").Value = "  +1.70%  "
$ws.Range("E26").Value = "  +2.81%  "
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.09"
$ws.Range("E28").Value = "  +0.83%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.59"
$ws.Range("E29").Value = "  +7.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "9.63"
$ws.Range("E30").Value = "  +1.66%  "
$ws.Range("E31").Value = "  -1.15%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "162.02"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("E33").Value = "  +3.85%  "
$ws.Range("E34").Value = "  +0.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0755"
$ws.Range("E35").Value = "  +1.96%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.10"
$ws.Range("E36").Value = "  +2.94%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "17.43"
$ws.Range("E37").Value = "  +5.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.108"
$ws.Range("E38").Value = "  +4.56%  "
$ws.Range("E39").Value = "  +5.54%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.39"
$ws.Range("E40").Value = "  +0.34%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.34"
$ws.Range("E43").Value = "  +15.21%  "
$ws.Range("D44").Value = "2.006.91"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "19.40"
$ws.Range("E45").Value = "  -0.65%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0285"
$ws.Range("E46").Value = "  +2.83%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.02"
$ws.Range("E47").Value = "  +6.09%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.24"
$ws.Range("E48").Value = "  +0.49%  "
$ws.Range("E49").Value = "  +5.55%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.56"
$ws.Range("E50").Value = "  +2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.83"
$ws.Range("E51").Value = "  +0.43%  "
